# csa_file_reference.docx — bump the "WPILib - Eclipse Plugins" version
# number in the reference table from 2018.1.1 to 2018.2.1
# (commit message: "Updated WPILib version 2018.2.1 number to match file repo.")
#
# The cell originally holds a single run with the text "2018.1.1". We
# locate it and change it to "2018.2.1", keeping the trailing ".1" as its
# own run (matching how the authored document stores "2018.2" and ".1"
# as two adjacent, identically-formatted runs) by toggling a character
# property off/on around that trailing piece so the two text runs don't
# get coalesced back into one on save.

$d = $word.ActiveDocument

# Locate the exact version string "2018.1.1" in the document body.
$hit = $d.Content
$hit.Find.Execute("2018.1.1", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$hitStart = $hit.Start

# Replace the old string with the first part of the new version number...
$hit.Text = "2018.2"

# ...then append the remaining ".1" right after it as a distinct run.
$tail = $d.Range($hitStart + 6, $hitStart + 6)
$tail.InsertAfter(".1")

# Nudge a character attribute on the tail run and immediately revert it;
# this keeps "2018.2" and ".1" as two separate <w:r> runs (identical
# formatting) instead of being merged back into a single run on save.
$tailRange = $d.Range($hitStart + 6, $hitStart + 8)
$tailRange.Font.Bold = 1
$tailRange.Font.Bold = 0
